{"js": "// Replace the division-problem text in each table cell with the new\n// values from the commit. Each source string is unique in the document,\n// so a matchCase/matchWholeWord search-and-replace per pair is safe and\n// unambiguous.\nconst replacements = [\n  [\"30\u00f73=\", \"53\u00f78=\"],\n  [\"94\u00f74=\", \"32\u00f72=\"],\n  [\"13\u00f74=\", \"67\u00f73=\"],\n  [\"95\u00f74=\", \"25\u00f77=\"],\n  [\"79\u00f77=\", \"55\u00f73=\"],\n  [\"70\u00f78=\", \"48\u00f76=\"],\n  [\"14\u00f72=\", \"62\u00f79=\"],\n  [\"17\u00f73=\", \"25\u00f77=\"],\n  [\"47\u00f77=\", \"71\u00f74=\"],\n  [\"26\u00f76=\", \"79\u00f76=\"],\n  [\"56\u00f78=\", \"50\u00f75=\"],\n  [\"60\u00f76=\", \"56\u00f73=\"],\n  [\"78\u00f77=\", \"75\u00f72=\"],\n  [\"83\u00f74=\", \"15\u00f79=\"],\n  [\"75\u00f77=\", \"36\u00f76=\"],\n  [\"88\u00f75=\", \"96\u00f72=\"],\n  [\"63\u00f76=\", \"83\u00f76=\"],\n  [\"62\u00f74=\", \"54\u00f77=\"],\n  [\"18\u00f79=\", \"78\u00f72=\"],\n  [\"27\u00f75=\", \"49\u00f78=\"],\n  [\"84\u00f78=\", \"35\u00f79=\"],\n  [\"13\u00f75=\", \"78\u00f79=\"],\n  [\"25\u00f72=\", \"71\u00f72=\"],\n  [\"56\u00f76=\", \"81\u00f76=\"],\n  [\"30\u00f75=\", \"77\u00f74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [from, to] of replacements) {\n  const results = body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(to, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the division-problem text in each table cell with the new\n# values from the commit. Each source string is unique in the document,\n# so a Find/Replace per pair is safe and unambiguous.\n\n$wdReplaceAll   = 2\n$wdFindContinue = 1\n\n$replacements = @(\n    @(\"30\u00f73=\", \"53\u00f78=\"),\n    @(\"94\u00f74=\", \"32\u00f72=\"),\n    @(\"13\u00f74=\", \"67\u00f73=\"),\n    @(\"95\u00f74=\", \"25\u00f77=\"),\n    @(\"79\u00f77=\", \"55\u00f73=\"),\n    @(\"70\u00f78=\", \"48\u00f76=\"),\n    @(\"14\u00f72=\", \"62\u00f79=\"),\n    @(\"17\u00f73=\", \"25\u00f77=\"),\n    @(\"47\u00f77=\", \"71\u00f74=\"),\n    @(\"26\u00f76=\", \"79\u00f76=\"),\n    @(\"56\u00f78=\", \"50\u00f75=\"),\n    @(\"60\u00f76=\", \"56\u00f73=\"),\n    @(\"78\u00f77=\", \"75\u00f72=\"),\n    @(\"83\u00f74=\", \"15\u00f79=\"),\n    @(\"75\u00f77=\", \"36\u00f76=\"),\n    @(\"88\u00f75=\", \"96\u00f72=\"),\n    @(\"63\u00f76=\", \"83\u00f76=\"),\n    @(\"62\u00f74=\", \"54\u00f77=\"),\n    @(\"18\u00f79=\", \"78\u00f72=\"),\n    @(\"27\u00f75=\", \"49\u00f78=\"),\n    @(\"84\u00f78=\", \"35\u00f79=\"),\n    @(\"13\u00f75=\", \"78\u00f79=\"),\n    @(\"25\u00f72=\", \"71\u00f72=\"),\n    @(\"56\u00f76=\", \"81\u00f76=\"),\n    @(\"30\u00f75=\", \"77\u00f74=\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $from = $pair[0]\n    $to   = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.Execute($from, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $to, $wdReplaceAll) | Out-Null\n}\n"}
